$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids numeric
# auto-coercion of numeric-looking strings like "598.31") and while
# preserving the cells original style (no left-over style change).
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Copy() | Out-Null
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.PasteSpecial(-4122) | Out-Null
}

Set-TextValue "D2" "68.409.23"
Set-TextValue "E2" "  +0.72%  "

Set-TextValue "D3" "2.645.68"
Set-TextValue "E3" "  +0.27%  "

Set-TextValue "E4" "  +0.01%  "

Set-TextValue "D5" "598.31"
Set-TextValue "E5" "  +0.07%  "

Set-TextValue "D6" "154.78"
Set-TextValue "E6" "  +0.70%  "

Set-TextValue "E7" "  -0.02%  "

Set-TextValue "E8" "  -0.68%  "

Set-TextValue "D9" "2.645.14"
Set-TextValue "E9" "  +0.29%  "

Set-TextValue "E10" "  +7.40%  "

Set-TextValue "E11" "  -0.54%  "

Set-TextValue "E12" "  +0.92%  "

Set-TextValue "E13" "  +2.10%  "

Set-TextValue "D14" "28.14"
Set-TextValue "E14" "  +1.64%  "

Set-TextValue "E15" "  +2.43%  "

Set-TextValue "D16" "3.128.56"
Set-TextValue "E16" "  +0.29%  "

Set-TextValue "D17" "68.271.48"
Set-TextValue "E17" "  +0.65%  "

Set-TextValue "D18" "2.645.69"
Set-TextValue "E18" "  +0.20%  "

Set-TextValue "D19" "11.37"
Set-TextValue "E19" "  -0.72%  "

Set-TextValue "D20" "363.85"
Set-TextValue "E20" "  -2.47%  "

Set-TextValue "D21" "7.49"
Set-TextValue "E21" "  -0.10%  "

Set-TextValue "D22" "4.38"
Set-TextValue "E22" "  +3.07%  "

Set-TextValue "E23" "  +2.09%  "

Set-TextValue "E24" "  +0.81%  "

Set-TextValue "D25" "74.68"
Set-TextValue "E25" "  +3.39%  "

Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  -0.04%  "

Set-TextValue "D27" "9.82"
Set-TextValue "E27" "  -1.24%  "

Set-TextValue "E28" "  +1.92%  "

Set-TextValue "D29" "2.776.95"
Set-TextValue "E29" "  +0.32%  "

Set-TextValue "E30" "  +0.14%  "

Set-TextValue "D31" "573.62"
Set-TextValue "E31" "  -1.12%  "

Set-TextValue "D32" "8.10"
Set-TextValue "E32" "  +2.93%  "

Set-TextValue "D33" "1.42"
Set-TextValue "E33" "  +1.45%  "

Set-TextValue "E34" "  +2.10%  "

Set-TextValue "E35" "  +3.05%  "

Set-TextValue "D36" "1.00"
Set-TextValue "E36" "  +0.02%  "

Set-TextValue "E37" "  +5.07%  "

Set-TextValue "D38" "160.97"
Set-TextValue "E38" "  +1.88%  "

Set-TextValue "D39" "19.36"
Set-TextValue "E39" "  +0.86%  "

Set-TextValue "E40" "  +1.64%  "

Set-TextValue "E41" "  -0.34%  "

Set-TextValue "D42" "5.37"
Set-TextValue "E42" "  +0.09%  "

Set-TextValue "B43" "dogwifhat"
Set-TextValue "C43" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "2.66"
Set-TextValue "E43" "  +1.13%  "

Set-TextValue "B44" "BabyDogeCoin"
Set-TextValue "C44" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D44" "0.0₆0335"
Set-TextValue "E44" "  -1.58%  "

Set-TextValue "E45" "  +3.33%  "

Set-TextValue "B46" "USDe"
Set-TextValue "C46" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D46" "1.00"
Set-TextValue "E46" "  +0.02%  "

Set-TextValue "B47" "OKB"
Set-TextValue "C47" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D47" "40.62"
Set-TextValue "E47" "  +0.93%  "

Set-TextValue "D48" "157.17"

Set-TextValue "D49" "3.76"
Set-TextValue "E49" "  +1.77%  "

Set-TextValue "E50" "  +0.70%  "

Set-TextValue "E51" "  +1.48%  "

$excel.CutCopyMode = 0
